# "adding averages and more checks"
# - Consolidate the bold title/header font into a single bold+white font
#   (title loses its 14pt size, header row text becomes white).
# - Training Dashboard: refresh PERIOD TO EXPIRE (H) / LAST UPDATE (I) for
#   rows 3-32 to reflect a check run 8 days later (16-Sep-2025 vs 08-Sep-2025).
# - Exam Dashboard: mark the dated rows as "date is valid" and narrow the
#   COMMENTS column.

$wb = $excel.ActiveWorkbook

# --- Shared look: bold white title + header row on both dashboards -------
foreach ($sheetName in @("Training Dashboard", "Exam Dashboard")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $titleCell = $ws.Cells.Item(1, 1)
    $titleCell.Font.Size = 11
    $titleCell.Font.Color = 16777215

    $lastCol = $ws.UsedRange.Columns.Count
    $headerRange = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item(2, $lastCol))
    $headerRange.Font.Color = 16777215
}

# --- Training Dashboard: refresh the expiry countdown --------------------
$training = $wb.Worksheets.Item("Training Dashboard")
for ($r = 3; $r -le 32; $r++) {
    $periodCell = $training.Cells.Item($r, 8)
    $periodCell.Value2 = $periodCell.Value2 - 8
    # Leading apostrophe keeps this a literal text value (matches the
    # original inline string) instead of Excel auto-parsing it as a date.
    $training.Cells.Item($r, 9).Value2 = "'16-Sep-2025"
}

# --- Exam Dashboard: new verdicts + narrower comments column --------------
$exam = $wb.Worksheets.Item("Exam Dashboard")
$exam.Range("E3").Value2 = "date is valid"
$exam.Range("E4").Value2 = "date is valid"
$exam.Range("E5").Value2 = "date is valid"
# ColumnWidth round-trips through a pixel conversion on save that adds
# ~5/6 of a character; pre-compensate so the saved width lands on 15.
$exam.Columns.Item(5).ColumnWidth = 15 - 5/6

Write-Host "edits applied"
